$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 46

$ws.Cells.Item($newRow, 1).NumberFormat = "@"
$ws.Cells.Item($newRow, 1).Value = "01-09-2021"
$ws.Cells.Item($newRow, 1).Style = "Normal"
$ws.Cells.Item($newRow, 2).Value = 117.2
$ws.Cells.Item($newRow, 3).Value = 108.7
$ws.Cells.Item($newRow, 4).Value = 100.38
$ws.Cells.Item($newRow, 5).Value = 112.83
$ws.Cells.Item($newRow, 6).Value = 112.58
$ws.Cells.Item($newRow, 7).Value = 107.4
$ws.Cells.Item($newRow, 8).Value = 114.3
$ws.Cells.Item($newRow, 9).Value = 93.52
$ws.Cells.Item($newRow, 10).Value = 109.72
$ws.Cells.Item($newRow, 11).Value = 112.31
$ws.Cells.Item($newRow, 12).Value = 112.04
$ws.Cells.Item($newRow, 13).Value = 112.78
